$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Clear deleted values in row 2 (B2:E2)
$ws.Range("B2:E2").ClearContents()

# Row 3 updates
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 2.8653260498306459
$ws.Range("D3").Value = -10.160790160633542
$ws.Range("E3").Value = 25.801246793190526

# Update selection to match new range
$ws.Range("B1:E3").Select()
